# [ADD] csv to excel converter
# Appends the rows produced by the CSV -> Excel conversion to the "users"
# sheet. The source CSV apparently got pasted in twice: once missing the
# leading "Date" column (rows 3-10, data lands in A:C) and once complete
# (rows 11-18, data lands in A:D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-10: Name, Country, Whatsapp (no Date column)
$shifted = @(
    @("roberto Santiago", "Afghanistan", "+93910966393"),
    @("Tatiana",           "Afghanistan", "+93919059173"),
    @("roberto Santiago", "Afghanistan", "+93910966392"),
    @("padre",             "Afghanistan", "+93654987245"),
    @("sdds",              "Afghanistan", "+9332323242342"),
    @("sdkn",              "Andorra",     "+37632323242342"),
    @("teste",             "Andorra",     "+376323232423411"),
    @("Thiciana Rocha ",   "Portugal",    "+351913895289")
)

$r = 3
foreach ($row in $shifted) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    # Leading "'" forces text so the "+" country-code prefix survives
    # instead of Excel coercing the value to a number.
    $ws.Cells.Item($r, 3).Value = "'" + $row[2]
    $r++
}

# Rows 11-18: Date, Name, Country, Whatsapp (full row, same source data)
$full = @(
    @("2024-07-11 03:40:15", "roberto Santiago", "Afghanistan", "+93910966393"),
    @("2024-07-11 03:40:15", "Tatiana",           "Afghanistan", "+93919059173"),
    @("2024-07-11 03:40:15", "roberto Santiago", "Afghanistan", "+93910966392"),
    @("2024-07-11 03:40:15", "padre",             "Afghanistan", "+93654987245"),
    @("2024-07-11 03:40:15", "sdds",              "Afghanistan", "+9332323242342"),
    @("2024-07-11 03:40:15", "sdkn",              "Andorra",     "+37632323242342"),
    @("2024-07-11 03:40:15", "teste",             "Andorra",     "+376323232423411"),
    @("2024-07-11 03:40:15", "Thiciana Rocha ",   "Portugal",    "+351913895289")
)

foreach ($row in $full) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    $r++
}
